# Apply the edits described by the diff:
# 1) Update several odds values in rows 2, 3 and 5.
# 2) Remove the three "ECUADOR - LIGA PRO" matches (rows 7-9), which shifts the
#    "USA - MLS" match (old row 10) up to become the new row 7.
# 3) Tweak two more values (M7, O7) on the resulting row 7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (YgMV5eZO / Tigre vs Instituto) updates ---
$ws.Range("G2").Value = 2.45
$ws.Range("I2").Value = 3.3
$ws.Range("J2").Value = 3.25
$ws.Range("W2").Value = 6.5
$ws.Range("AH2").Value = 7.5
$ws.Range("AL2").Value = 29
$ws.Range("AS2").Value = 301
$ws.Range("AX2").Value = 19

# --- Row 3 (vZ5Qsl5t / Tomayapo vs SA Bulo Bulo) updates ---
$ws.Range("G3").Value = 2
$ws.Range("H3").Value = 3.5
$ws.Range("I3").Value = 3.6
$ws.Range("J3").Value = 2.63
$ws.Range("L3").Value = 4
$ws.Range("N3").Value = 12
$ws.Range("X3").Value = 10
$ws.Range("Z3").Value = 17
$ws.Range("AA3").Value = 15
$ws.Range("AJ3").Value = 13
$ws.Range("AK3").Value = 41
$ws.Range("AL3").Value = 29
$ws.Range("AN3").Value = 4
$ws.Range("AQ3").Value = 34
$ws.Range("BA3").Value = 81

# --- Row 5 (IcQ1l8Is / Botafogo RJ vs Vitoria) updates ---
$ws.Range("H5").Value = 5
$ws.Range("I5").Value = 9
$ws.Range("N5").Value = 12
$ws.Range("Y5").Value = 9
$ws.Range("AF5").Value = 81

# --- Remove the three ECUADOR - LIGA PRO rows (7, 8, 9). ---
# Deleting the same row index three times removes the original rows 7, 8 and 9
# and shifts the remaining row (old row 10, USA - MLS) up into row 7.
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(7).Delete()

# --- Final tweaks on the new row 7 (previously row 10) ---
$ws.Range("M7").Value = 1.03
$ws.Range("O7").Value = 1.25
